$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D:E) store values as plain text (e.g. "26.370.89",
# "  -2.75%  ") using European-style "." thousand separators and percent
# strings with padding spaces. Excel's COM Range.Value setter auto-detects
# numeric-looking strings and silently converts them to real numbers, which
# would corrupt values like "308.34" or "1.002" (dropping the literal text
# form / float precision) and strip the leading/trailing spaces + "%" off
# percentage cells that happen to look numeric once trimmed.
#
# Force the data rows (D2:E51) to Text format first so every write below is
# stored as a literal string, matching the original file's inlineStr cells.
# Row 1 (headers) and column A (index numbers) are intentionally excluded so
# their existing styles are left untouched.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.350.88"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").Value = "1.774.09"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "308.34"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("D8").Value = "0.3607"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "0.07127"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "0.8371"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "20.37"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.768.04"
$ws.Range("E12").Value = "  -5.98%  "
$ws.Range("D13").Value = "6.446"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "5.247"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "78.97"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "0.000008652"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "26.370.31"
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").Value = "5.098"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "1.990.32"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").Value = "151.84"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("E26").Value = "  -8.10%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "5.064"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "114.25"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "1.841"
$ws.Range("E30").Value = "  +11.34%  "
$ws.Range("D31").Value = "0.08832"
$ws.Range("D32").Value = "0.7260"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "1.120"
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("D34").Value = "4.311"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "2.737"
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "0.01887"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.1608"
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.4910"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "2.594"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "6.330"
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("D44").Value = "8.038"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "104.68"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "1.623"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("D49").Value = "0.06172"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "1.696"
$ws.Range("E51").Value = "  +1.88%  "

# Restore the default "Normal" style on the touched range so no stray
# number-format style lingers on these cells (they had no style index
# originally).
$dataRange.Style = "Normal"
